$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text to avoid locale numeric parsing
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.883.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.873.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.887.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.249.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.536"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.835"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.808"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.786.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("D51").Style = "Normal"

# Volume (column E) updates
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +7.23%  "
$ws.Range("E23").Value = "  +6.17%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("E47").Value = "  +1.09%  "
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("E51").Value = "  +0.69%  "
